$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove all existing hyperlinks first so we can rebuild them cleanly
$ws.Hyperlinks.Delete()

# Row 2
$ws.Cells.Item(2,1).Value = "2026-01-21 18:39:42"
$ws.Cells.Item(2,2).Value = "大手SIer等のAIソリューション開発・導入を支援してくださるエンジニア・PM募集"
$ws.Cells.Item(2,3).Value = "システム開発"
$ws.Cells.Item(2,4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(2,5).Value = "期限情報なし"
$ws.Cells.Item(2,6).Value = "https://www.lancers.jp/work/detail/5455098"
$ws.Cells.Item(2,7).Value = 375
$ws.Cells.Item(2,8).Value = "🔥AI,Ai ◆開発"

# Row 3
$ws.Cells.Item(3,1).Value = "2026-01-21 18:39:42"
$ws.Cells.Item(3,2).Value = "法人向け生成AIサービス(RAG・議事録機能)の設計・開発を支援エンジニア募集(AI/バックエンド)"
$ws.Cells.Item(3,3).Value = "システム開発"
$ws.Cells.Item(3,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(3,5).Value = "期限情報なし"
$ws.Cells.Item(3,6).Value = "https://www.lancers.jp/work/detail/5445159"
$ws.Cells.Item(3,7).Value = 368
$ws.Cells.Item(3,8).Value = "🔥AI,Ai ◆開発"

# Row 4
$ws.Cells.Item(4,1).Value = "2026-01-21 18:39:42"
$ws.Cells.Item(4,2).Value = "B2B向け生成AIサービス(チャット・RAG)の新規開発プロジェクト推進を支援してくださるPM募集"
$ws.Cells.Item(4,3).Value = "システム開発"
$ws.Cells.Item(4,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(4,5).Value = "期限情報なし"
$ws.Cells.Item(4,6).Value = "https://www.lancers.jp/work/detail/5445154"
$ws.Cells.Item(4,7).Value = 368
$ws.Cells.Item(4,8).Value = "🔥AI,Ai ◆開発"

# Row 5
$ws.Cells.Item(5,1).Value = "2026-01-21 18:39:42"
$ws.Cells.Item(5,2).Value = "【急募】マッチングアプリのLLMO・AIO継続支援をお手伝いください!"
$ws.Cells.Item(5,3).Value = "システム開発"
$ws.Cells.Item(5,4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(5,5).Value = "期限情報なし"
$ws.Cells.Item(5,6).Value = "https://www.lancers.jp/work/detail/5476280"
$ws.Cells.Item(5,7).Value = 333
$ws.Cells.Item(5,8).Value = "🔥AI,Ai ◇アプリ"

# Row 6
$ws.Cells.Item(6,1).Value = "2026-01-21 18:39:42"
$ws.Cells.Item(6,2).Value = "マッチングアプリのLLMO・AIO対策を継続支援いただける方を募集"
$ws.Cells.Item(6,3).Value = "システム開発"
$ws.Cells.Item(6,4).Value = "20,000 円 ~"
$ws.Cells.Item(6,5).Value = "期限情報なし"
$ws.Cells.Item(6,6).Value = "https://www.lancers.jp/work/detail/5476284"
$ws.Cells.Item(6,7).Value = 330
$ws.Cells.Item(6,8).Value = "🔥AI,Ai ◇アプリ"

# Row 7
$ws.Cells.Item(7,1).Value = "2026-01-21 18:39:42"
$ws.Cells.Item(7,2).Value = "【長期案件】生成AIを利用したチャットボット作成のPMOを募集"
$ws.Cells.Item(7,3).Value = "システム開発"
$ws.Cells.Item(7,4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(7,5).Value = "期限情報なし"
$ws.Cells.Item(7,6).Value = "https://www.lancers.jp/work/detail/5476159"
$ws.Cells.Item(7,7).Value = 310
$ws.Cells.Item(7,8).Value = "🔥AI,Ai"

# Row 8
$ws.Cells.Item(8,1).Value = "2026-01-21 18:39:42"
$ws.Cells.Item(8,2).Value = "【急募】フットアールサッカースクール向け出欠管理Webアプリ開発者募集"
$ws.Cells.Item(8,3).Value = "システム開発"
$ws.Cells.Item(8,4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(8,5).Value = "期限情報なし"
$ws.Cells.Item(8,6).Value = "https://www.lancers.jp/work/detail/5475081"
$ws.Cells.Item(8,7).Value = 128
$ws.Cells.Item(8,8).Value = "◆開発 ◇アプリ"

# Row 9
$ws.Cells.Item(9,1).Value = "2026-01-21 18:39:42"
$ws.Cells.Item(9,2).Value = "【継続依頼あり】教育システム開発案件のクロージング代行&要件定義"
$ws.Cells.Item(9,3).Value = "システム開発"
$ws.Cells.Item(9,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(9,5).Value = "期限情報なし"
$ws.Cells.Item(9,6).Value = "https://www.lancers.jp/work/detail/5475937"
$ws.Cells.Item(9,7).Value = 118
$ws.Cells.Item(9,8).Value = "◆開発,システム開発"

# Row 10
$ws.Cells.Item(10,1).Value = "2026-01-21 18:39:42"
$ws.Cells.Item(10,2).Value = "【急募】野球スコアボードシステム開発のフリーランス募集"
$ws.Cells.Item(10,3).Value = "システム開発"
$ws.Cells.Item(10,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(10,5).Value = "期限情報なし"
$ws.Cells.Item(10,6).Value = "https://www.lancers.jp/work/detail/5475665"
$ws.Cells.Item(10,7).Value = 118
$ws.Cells.Item(10,8).Value = "◆開発,システム開発"

# Row 11
$ws.Cells.Item(11,1).Value = "2026-01-21 18:39:42"
$ws.Cells.Item(11,2).Value = "【急募】自動車整備業向けCRM構築パートナー募集"
$ws.Cells.Item(11,3).Value = "システム開発"
$ws.Cells.Item(11,4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(11,5).Value = "期限情報なし"
$ws.Cells.Item(11,6).Value = "https://www.lancers.jp/work/detail/5474125"
$ws.Cells.Item(11,7).Value = 25
$ws.Cells.Item(11,8).Value = ""

# Row 12
$ws.Cells.Item(12,1).Value = "2026-01-21 18:39:42"
$ws.Cells.Item(12,2).Value = "【急募】CSVデータをワードに自動入力するスキルをお持ちの方"
$ws.Cells.Item(12,3).Value = "システム開発"
$ws.Cells.Item(12,4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(12,5).Value = "期限情報なし"
$ws.Cells.Item(12,6).Value = "https://www.lancers.jp/work/detail/5475924"
$ws.Cells.Item(12,7).Value = 10
$ws.Cells.Item(12,8).Value = ""

# Re-create hyperlinks for column F (rows 2-12)
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5455098") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5445159") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5445154") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5476280") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5476284") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5476159") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5475081") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5475937") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5475665") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5474125") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5475924") | Out-Null

# Column H width: 12 -> 13
$ws.Range("H1").ColumnWidth = 12.14
